$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 147; this pushes the existing rows 147..212 down to 148..213
# (mirrors a new weekly record being prepended to this variety's price history).
$ws.Rows("147:147").Insert()

# Populate the newly inserted row 147 with the new weekly record.
# The "dimension"/header-like columns (A,B,C,E,F,G,H,I,J,K,L,R) keep the same
# constant values used throughout this block; D is the new date, and
# M,N,O,P,Q,S,T carry the values that previously lived in (old) row 146.
$ws.Range("A147").Value = 5
$ws.Range("B147").Value = "Macroferia Regional de Talca"
$ws.Range("C147").Value = "Maule"
$ws.Range("D147").Value = 44609
$ws.Range("D147").NumberFormat = $ws.Range("D148").NumberFormat
$ws.Range("E147").Value = 7
$ws.Range("F147").Value = "Fruta"
$ws.Range("G147").Value = 100108
$ws.Range("H147").Value = "Tropicales y subtropicales"
$ws.Range("I147").Value = 100108005
$ws.Range("J147").Value = "Piña"
$ws.Range("K147").Value = "Caramelo"
$ws.Range("L147").Value = "Segunda"
$ws.Range("M147").Value = 300
$ws.Range("N147").Value = 15000
$ws.Range("O147").Value = 15000
$ws.Range("P147").Value = 15000
$ws.Range("Q147").Value = "$/caja 14 unidades"
$ws.Range("R147").Value = "Ecuador"
$ws.Range("S147").Value = 1071
$ws.Range("T147").Value = 14
